$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting already used in column A (A2) down onto the new rows
# A3:A5 before changing the values, so the whole column keeps consistent
# styling (bold, border, centered) like the original A2/A3 cells.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# Update existing B2 value (was 12 -> 6)
$ws.Range("B2").Value = 6

# Row 3 now holds new values (3,3) instead of the old (1,1)
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 3

# New row 4 (4,3)
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 3

# Old row 3 data (1,1) moves down to row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 1
